$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had two columns:
#   A: "button_apiKeyManagement_trNthChild" / "2"
#   B: "input_KeyName"                     / (blank)
# Column B is removed; its header becomes the (sole remaining) column A
# header, and row 2 is left blank.

$newHeader = $ws.Range("B1").Value2
$ws.Range("A1").Value = $newHeader

# Row 2 becomes blank.
$ws.Range("A2").Value = ""
# Touch formatting on the now-blank cell so it stays present in the sheet
# (rather than being dropped) - this mirrors the blank placeholder cell
# that used to live at B2.
$ws.Range("A2").Font.Bold = $false

# Drop column B's contents/formatting entirely - only column A remains.
$ws.Columns("B").Clear()

# Column A keeps the narrower width that used to belong to column B.
$ws.Columns("A").ColumnWidth = 14.17
